$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20. This pushes the existing
# "Total" row (formerly row 20) down to row 21.
$ws.Rows("20:20").Insert()

# Copy the formatting (styles) of the last data row (19) onto the
# newly inserted row 20 so it matches the other data rows.
$ws.Range("A19:F19").Copy()
$ws.Range("A20:F20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new data row: 2024-01-01, 16:00 -> 21:00, rate 10.
$ws.Range("A20").Value = 45292
$ws.Range("B20").Value = 0.666666666666667
$ws.Range("C20").Value = 0.875
$ws.Range("D20").Formula = "=(C20<B20)+C20-B20"
$ws.Range("E20").Value = 10
$ws.Range("F20").Formula = "=(D20*24)*E20"

# Entering the F20 formula (which references the time-formatted D20
# cell) causes the engine to infer a time number format for F20;
# restore it to the "General" format used by the other Bill cells.
$ws.Range("F20").NumberFormat = "General"

# Update the Total row (now row 21) so its SUM ranges include the
# newly added row 20.
$ws.Range("D21").Formula = "=SUM(D2:D20)"
$ws.Range("F21").Formula = "=SUM(F2:F20)"

# Move the active selection to F22, matching where Excel would leave
# the cursor after the row was inserted above the total row.
$ws.Range("F22").Select() | Out-Null
